# Apply formatting updates to the Income Statement template:
#  - zoom to 125%, move active selection to K10
#  - left-align the "Row Index" / "Summary Index" header cells (A1/E1)
#  - left-align columns A and E for all data rows
#  - right-align columns B/C/D for regular (non-subtotal) rows, incl. row 31
#  - left-align columns B/C/D for the bold subtotal rows (5,12,14,15,28,30)

$xlHAlignLeft  = -4131
$xlHAlignRight = -4152

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-HAlign($range, $align) {
    # Multi-area ranges only apply the format to the first area in this
    # runtime, so walk the Areas collection explicitly.
    foreach ($area in $range.Areas) {
        $area.HorizontalAlignment = $align
    }
}

# --- Sheet view: zoom + selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 125
$ws.Range("K10").Select()

# --- Header row (row 1): A1/E1 become left-aligned (keep vertical=center), B1:D1 stay centered ---
Set-HAlign $ws.Range("A1,E1") $xlHAlignLeft

# --- Columns A and E (rows 2-31): left aligned ---
Set-HAlign $ws.Range("A2:A31,E2:E31") $xlHAlignLeft

# --- Columns B:D for normal rows (incl. row 31 "Net Profit %"): right aligned ---
Set-HAlign $ws.Range("B2:D4,B6:D11,B13:D13,B16:D27,B29:D29,B31:D31") $xlHAlignRight

# --- Columns B:D for bold subtotal rows (5,12,14,15,28,30): left aligned ---
Set-HAlign $ws.Range("B5:D5,B12:D12,B14:D15,B28:D28,B30:D30") $xlHAlignLeft
